# Adds the _06 and _07 series of model runs to the all_runs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# ---------------------------------------------------------------------------
# Helper to fill one data row (columns A-H) in one shot.
# Positional params: Row, Project, Year, Directory, RunSet, Category,
#                     UrbansimPath, UrbansimRunId, Status
# ---------------------------------------------------------------------------
function Set-RunRow {
    param($Row, $Project, $Year, $Directory, $RunSet, $Category, $UrbansimPath, $UrbansimRunId, $Status)
    $ws.Cells.Item($Row, 1).Value = $Project
    $ws.Cells.Item($Row, 2).Value = $Year
    $ws.Cells.Item($Row, 3).Value = $Directory
    $ws.Cells.Item($Row, 4).Value = $RunSet
    $ws.Cells.Item($Row, 5).Value = $Category
    $ws.Cells.Item($Row, 6).Value = $UrbansimPath
    $ws.Cells.Item($Row, 7).Value = $UrbansimRunId
    $ws.Cells.Item($Row, 8).Value = $Status
}

# ---------------------------------------------------------------------------
# 2035 No Project block: clear the stale "current"/"running" statuses on the
# existing _02 / _03 rows, then insert the new _04 / _06 / _07 rows right
# before the "Basic" rows.
# ---------------------------------------------------------------------------
$ws.Cells.Item(25, 8).Value = ""
$ws.Cells.Item(26, 8).Value = ""

$ws.Rows(27).Insert()
$ws.Rows(27).Insert()
$ws.Rows(27).Insert()

Set-RunRow 27 "RTP2021" 2035 "2035_TM152_DBP_NoProject_04" "DraftBlueprint" "No Project" '"Blueprint Plus Crossing (s23)\v1.5.5"' "run998" ""
Set-RunRow 28 "RTP2021" 2035 "2035_TM152_DBP_NoProject_06" "DraftBlueprint" "No Project" '"Blueprint Plus Crossing (s23)\v1.6 (all strategies)"' "run90" "current"
Set-RunRow 29 "RTP2021" 2035 "2035_TM152_DBP_NoProject_07" "DraftBlueprint" "No Project" '"Blueprint Plus Crossing (s23)\v1.7 (strategies + BASIS-hybrid)"' "run92" ""

# ---------------------------------------------------------------------------
# 2035 Plus block: insert the new _06 / _07 rows right after _04, before the
# 2050 section begins.
# ---------------------------------------------------------------------------
$ws.Rows(37).Insert()
$ws.Rows(37).Insert()

Set-RunRow 37 "RTP2021" 2035 "2035_TM152_DBP_Plus_06" "DraftBlueprint" "Plus" '"Blueprint Plus Crossing (s23)\v1.6 (all strategies)"' "run90" "current"
Set-RunRow 38 "RTP2021" 2035 "2035_TM152_DBP_Plus_07" "DraftBlueprint" "Plus" '"Blueprint Plus Crossing (s23)\v1.7 (strategies + BASIS-hybrid)"' "run92" ""

# ---------------------------------------------------------------------------
# 2050 No Project block: insert the new _06 / _07 rows right before the 2050
# "Basic" rows.
# ---------------------------------------------------------------------------
$ws.Rows(43).Insert()
$ws.Rows(43).Insert()

Set-RunRow 43 "RTP2021" 2050 "2050_TM152_DBP_NoProject_06" "DraftBlueprint" "No Project" '"Blueprint Plus Crossing (s23)\v1.6 (all strategies)"' "run90" "current"
Set-RunRow 44 "RTP2021" 2050 "2050_TM152_DBP_NoProject_07" "DraftBlueprint" "No Project" '"Blueprint Plus Crossing (s23)\v1.7 (strategies + BASIS-hybrid)"' "run92" ""

# ---------------------------------------------------------------------------
# 2050 Plus Crossing block: append the new _06 / _07 rows at the very end.
# Row 54 (RTP2022 / 2051) is a new, separate project entry.
# ---------------------------------------------------------------------------
$ws.Rows(53).Insert()
$ws.Rows(53).Insert()

Set-RunRow 53 "RTP2021" 2050 "2050_TM152_DBP_PlusCrossing_06" "DraftBlueprint" "Plus" '"Blueprint Plus Crossing (s23)\v1.6 (all strategies)"' "run90" "current"
Set-RunRow 54 "RTP2022" 2051 "2050_TM152_DBP_PlusCrossing_07" "DraftBlueprint" "Plus" '"Blueprint Plus Crossing (s23)\v1.7 (strategies + BASIS-hybrid)"' "run92" ""

# ---------------------------------------------------------------------------
# Refresh the view state to match (frozen pane scrolled down one row, active
# cell tracking the bottom of the newly expanded table).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A24").Select()
$ws.Range("H55").Select()
